# Auto-generated cell updates applying the diff to cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''30.101.68'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.56%  '
$ws.Range("D3").Value = '''1.641.54'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.39%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''215.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.49%  '
$ws.Range("E6").Value = '  +1.39%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '''29.31'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.55%  '
$ws.Range("E9").Value = '  +4.03%  '
$ws.Range("E10").Value = '  +2.23%  '
$ws.Range("D11").Value = '''0.0917'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.67%  '
$ws.Range("D12").Value = '''1.876.42'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.43%  '
$ws.Range("D13").Value = '''1.648.77'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.89%  '
$ws.Range("D14").Value = '''0.575'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.84%  '
$ws.Range("D15").Value = '''9.48'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +23.59%  '
$ws.Range("E16").Value = '  +4.29%  '
$ws.Range("D17").Value = '''30.118.90'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.68%  '
$ws.Range("D18").Value = '''65.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.04%  '
$ws.Range("D19").Value = '''248.25'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.86%  '
$ws.Range("D20").Value = '''0.0₃0711'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.24%  '
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("D22").Value = '''4.20'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.88%  '
$ws.Range("D23").Value = '''9.94'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.64%  '
$ws.Range("E24").Value = '  +1.58%  '
$ws.Range("D25").Value = '''159.34'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.81%  '
$ws.Range("D26").Value = '''15.77'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.34%  '
$ws.Range("E27").Value = '  +2.93%  '
$ws.Range("D28").Value = '''6.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.22%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  +2.71%  '
$ws.Range("E31").Value = '  +5.99%  '
$ws.Range("E32").Value = '  +5.78%  '
$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("D34").Value = '''1.438.98'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.20%  '
$ws.Range("E35").Value = '  +7.38%  '
$ws.Range("E36").Value = '  +1.67%  '
$ws.Range("E37").Value = '  -0.21%  '
$ws.Range("D38").Value = '''77.53'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +17.40%  '
$ws.Range("E39").Value = '  +1.58%  '
$ws.Range("E40").Value = '  +0.07%  '
$ws.Range("E41").Value = '  +2.95%  '
$ws.Range("D42").Value = '''2.04'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.11%  '
$ws.Range("E43").Value = '  +3.15%  '
$ws.Range("D44").Value = '''55.57'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.75%  '
$ws.Range("E45").Value = '  +1.01%  '
$ws.Range("E46").Value = '  +5.68%  '
$ws.Range("D47").Value = '''1.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("E48").Value = '  +1.46%  '
$ws.Range("D49").Value = '''1.783.17'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.45%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '''0.0₆0114'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.55%  '
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").Value = '''90.44'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.67%  '
